$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{Cell="B2"; Value=1.02},
    @{Cell="C2"; Value=1.03901134014721},
    @{Cell="D2"; Value=1.040098293597874},
    @{Cell="E2"; Value=1.047235029935547},
    @{Cell="F2"; Value=1.057362519252723},
    @{Cell="I2"; Value=1.038234503270647},
    @{Cell="J2"; Value=1.044105668756326},
    @{Cell="K2"; Value=1.042881606304943},
    @{Cell="L2"; Value=1.049998235410718},
    @{Cell="M2"; Value=1.060097696066628},
    @{Cell="N2"; Value=1.018590181221952},
    @{Cell="B3"; Value=1.02},
    @{Cell="C3"; Value=1.039989348756141},
    @{Cell="D3"; Value=1.040811860595336},
    @{Cell="E3"; Value=1.048177832135185},
    @{Cell="F3"; Value=1.058555236205481},
    @{Cell="I3"; Value=1.038465558715381},
    @{Cell="J3"; Value=1.044728564242412},
    @{Cell="K3"; Value=1.043405782001332},
    @{Cell="L3"; Value=1.050752501023985},
    @{Cell="M3"; Value=1.06110327184548},
    @{Cell="N3"; Value=1.018800264955289},
    @{Cell="B4"; Value=1.02},
    @{Cell="C4"; Value=1.040622402358694},
    @{Cell="D4"; Value=1.041273723083823},
    @{Cell="E4"; Value=1.048788495110687},
    @{Cell="F4"; Value=1.059328111547067},
    @{Cell="I4"; Value=1.038613975361191},
    @{Cell="J4"; Value=1.045131227876188},
    @{Cell="K4"; Value=1.04374442382901},
    @{Cell="L4"; Value=1.05124054347213},
    @{Cell="M4"; Value=1.061754478872474},
    @{Cell="N4"; Value=1.018935987745562},
    @{Cell="B5"; Value=1.02},
    @{Cell="C5"; Value=1.040888589530861},
    @{Cell="D5"; Value=1.041467921886013},
    @{Cell="E5"; Value=1.049045362468875},
    @{Cell="F5"; Value=1.059653293258789},
    @{Cell="I5"; Value=1.038676108031265},
    @{Cell="J5"; Value=1.04530041342232},
    @{Cell="K5"; Value=1.043886660191496},
    @{Cell="L5"; Value=1.051445711744768},
    @{Cell="M5"; Value=1.062028373378309},
    @{Cell="N5"; Value=1.018992993769565},
    @{Cell="B6"; Value=1.02},
    @{Cell="C6"; Value=1.040933286532425},
    @{Cell="D6"; Value=1.041500530556733},
    @{Cell="E6"; Value=1.049088500091657},
    @{Cell="F6"; Value=1.059707908226226},
    @{Cell="I6"; Value=1.038686525024682},
    @{Cell="J6"; Value=1.045328814877598},
    @{Cell="K6"; Value=1.043910534740243},
    @{Cell="L6"; Value=1.051480160106149},
    @{Cell="M6"; Value=1.062074368912474},
    @{Cell="N6"; Value=1.019002562284476},
    @{Cell="B7"; Value=1.02},
    @{Cell="C7"; Value=1.040625958963121},
    @{Cell="D7"; Value=1.041276317853937},
    @{Cell="E7"; Value=1.048791926817358},
    @{Cell="F7"; Value=1.059332455600392},
    @{Cell="I7"; Value=1.038614806609641},
    @{Cell="J7"; Value=1.045133488912732},
    @{Cell="K7"; Value=1.043746324904557},
    @{Cell="L7"; Value=1.051243284960005},
    @{Cell="M7"; Value=1.06175813816539},
    @{Cell="N7"; Value=1.018936749666857},
    @{Cell="B8"; Value=1.02},
    @{Cell="C8"; Value=1.039341817623197},
    @{Cell="D8"; Value=1.040339417972141},
    @{Cell="E8"; Value=1.047553528295965},
    @{Cell="F8"; Value=1.057765374005695},
    @{Cell="I8"; Value=1.038312815375221},
    @{Cell="J8"; Value=1.044316259916118},
    @{Cell="K8"; Value=1.043058864325409},
    @{Cell="L8"; Value=1.050253146397777},
    @{Cell="M8"; Value=1.060437424689061},
    @{Cell="N8"; Value=1.018661224452647},
    @{Cell="B9"; Value=1.02},
    @{Cell="C9"; Value=1.037080674253375},
    @{Cell="D9"; Value=1.038689571965601},
    @{Cell="E9"; Value=1.045375988695186},
    @{Cell="F9"; Value=1.055012471626614},
    @{Cell="I9"; Value=1.037772322187813},
    @{Cell="J9"; Value=1.042873224256397},
    @{Cell="K9"; Value=1.041843400903896},
    @{Cell="L9"; Value=1.048508281114122},
    @{Cell="M9"; Value=1.058114250613558},
    @{Cell="N9"; Value=1.018174075363549},
    @{Cell="B10"; Value=1.02},
    @{Cell="C10"; Value=1.035574396223345},
    @{Cell="D10"; Value=1.037590464670735},
    @{Cell="E10"; Value=1.043927476887079},
    @{Cell="F10"; Value=1.053182928724752},
    @{Cell="I10"; Value=1.037406401426043},
    @{Cell="J10"; Value=1.041909231949815},
    @{Cell="K10"; Value=1.04103039031022},
    @{Cell="L10"; Value=1.047344988592202},
    @{Cell="M10"; Value=1.056568237533042},
    @{Cell="N10"; Value=1.017848222761054},
    @{Cell="B11"; Value=1.02},
    @{Cell="C11"; Value=1.034922438283594},
    @{Cell="D11"; Value=1.037114737983313},
    @{Cell="E11"; Value=1.043301017027008},
    @{Cell="F11"; Value=1.052392076415302},
    @{Cell="I11"; Value=1.037246631229231},
    @{Cell="J11"; Value=1.041491351784645},
    @{Cell="K11"; Value=1.040677713932848},
    @{Cell="L11"; Value=1.046841263858447},
    @{Cell="M11"; Value=1.055899456913695},
    @{Cell="N11"; Value=1.017706870303027},
    @{Cell="B12"; Value=1.02},
    @{Cell="C12"; Value=1.034680312939221},
    @{Cell="D12"; Value=1.036938061944609},
    @{Cell="E12"; Value=1.043068435930988},
    @{Cell="F12"; Value=1.052098521817473},
    @{Cell="I12"; Value=1.0371870868049},
    @{Cell="J12"; Value=1.041336063030177},
    @{Cell="K12"; Value=1.040546618967821},
    @{Cell="L12"; Value=1.046654156718533},
    @{Cell="M12"; Value=1.055651140327091},
    @{Cell="N12"; Value=1.017654327449139},
    @{Cell="B13"; Value=1.02},
    @{Cell="C13"; Value=1.034732247794468},
    @{Cell="D13"; Value=1.036975958181321},
    @{Cell="E13"; Value=1.043118320199102},
    @{Cell="F13"; Value=1.052161481080579},
    @{Cell="I13"; Value=1.037199868271443},
    @{Cell="J13"; Value=1.041369376131271},
    @{Cell="K13"; Value=1.04057474358514},
    @{Cell="L13"; Value=1.046694291890165},
    @{Cell="M13"; Value=1.055704400626203},
    @{Cell="N13"; Value=1.017665599797292},
    @{Cell="B14"; Value=1.02},
    @{Cell="C14"; Value=1.034902423292098},
    @{Cell="D14"; Value=1.037100133269025},
    @{Cell="E14"; Value=1.04328178947274},
    @{Cell="F14"; Value=1.052367806955822},
    @{Cell="I14"; Value=1.037241713319561},
    @{Cell="J14"; Value=1.041478516982145},
    @{Cell="K14"; Value=1.040666879522371},
    @{Cell="L14"; Value=1.046825797549207},
    @{Cell="M14"; Value=1.055878928980628},
    @{Cell="N14"; Value=1.017702527873477},
    @{Cell="B15"; Value=1.02},
    @{Cell="C15"; Value=1.035007279454689},
    @{Cell="D15"; Value=1.037176645645203},
    @{Cell="E15"; Value=1.043382523407432},
    @{Cell="F15"; Value=1.05249495806021},
    @{Cell="I15"; Value=1.037267469118588},
    @{Cell="J15"; Value=1.041545753071799},
    @{Cell="K15"; Value=1.040723634907173},
    @{Cell="L15"; Value=1.046906822358713},
    @{Cell="M15"; Value=1.05598647469254},
    @{Cell="N15"; Value=1.017725275417225},
    @{Cell="B16"; Value=1.02},
    @{Cell="C16"; Value=1.035617670241781},
    @{Cell="D16"; Value=1.037622041247621},
    @{Cell="E16"; Value=1.043969068921993},
    @{Cell="F16"; Value=1.053235443474953},
    @{Cell="I16"; Value=1.037416976973746},
    @{Cell="J16"; Value=1.041936955516598},
    @{Cell="K16"; Value=1.041053782893595},
    @{Cell="L16"; Value=1.047378418934939},
    @{Cell="M16"; Value=1.056612636077716},
    @{Cell="N16"; Value=1.017857598480854},
    @{Cell="B17"; Value=1.02},
    @{Cell="C17"; Value=1.036000624814853},
    @{Cell="D17"; Value=1.037901478712565},
    @{Cell="E17"; Value=1.044337196006101},
    @{Cell="F17"; Value=1.053700292381656},
    @{Cell="I17"; Value=1.037510404994521},
    @{Cell="J17"; Value=1.042182222245302},
    @{Cell="K17"; Value=1.041260705706256},
    @{Cell="L17"; Value=1.047674236214914},
    @{Cell="M17"; Value=1.057005585582888},
    @{Cell="N17"; Value=1.017940532827029},
    @{Cell="B18"; Value=1.02},
    @{Cell="C18"; Value=1.036224021810348},
    @{Cell="D18"; Value=1.038064488519409},
    @{Cell="E18"; Value=1.044551991148184},
    @{Cell="F18"; Value=1.053971561312694},
    @{Cell="I18"; Value=1.037564772204409},
    @{Cell="J18"; Value=1.042325237271094},
    @{Cell="K18"; Value=1.041381338745448},
    @{Cell="L18"; Value=1.047846780305286},
    @{Cell="M18"; Value=1.05723484939824},
    @{Cell="N18"; Value=1.017988882300115},
    @{Cell="B19"; Value=1.02},
    @{Cell="C19"; Value=1.036300198827749},
    @{Cell="D19"; Value=1.038120073789329},
    @{Cell="E19"; Value=1.044625243080575},
    @{Cell="F19"; Value=1.054064079129017},
    @{Cell="I19"; Value=1.037583288337131},
    @{Cell="J19"; Value=1.042373994069022},
    @{Cell="K19"; Value=1.041422461041191},
    @{Cell="L19"; Value=1.047905613179369},
    @{Cell="M19"; Value=1.057313033118273},
    @{Cell="N19"; Value=1.018005364019894},
    @{Cell="B20"; Value=1.02},
    @{Cell="C20"; Value=1.035959534714831},
    @{Cell="D20"; Value=1.037871495780493},
    @{Cell="E20"; Value=1.044297691915727},
    @{Cell="F20"; Value=1.053650405001873},
    @{Cell="I20"; Value=1.037500394266369},
    @{Cell="J20"; Value=1.042155912088353},
    @{Cell="K20"; Value=1.041238511208044},
    @{Cell="L20"; Value=1.047642497933933},
    @{Cell="M20"; Value=1.056963419300347},
    @{Cell="N20"; Value=1.017931637312659},
    @{Cell="B21"; Value=1.02},
    @{Cell="C21"; Value=1.034852309719591},
    @{Cell="D21"; Value=1.037063565964472},
    @{Cell="E21"; Value=1.043233648695234},
    @{Cell="F21"; Value=1.052307043519275},
    @{Cell="I21"; Value=1.037229396477642},
    @{Cell="J21"; Value=1.041446379635692},
    @{Cell="K21"; Value=1.0406397504107},
    @{Cell="L21"; Value=1.046787072444495},
    @{Cell="M21"; Value=1.055827532014948},
    @{Cell="N21"; Value=1.017691654529741},
    @{Cell="B22"; Value=1.02},
    @{Cell="C22"; Value=1.034156389961887},
    @{Cell="D22"; Value=1.036555762631648},
    @{Cell="E22"; Value=1.042565302252304},
    @{Cell="F22"; Value=1.051463594309597},
    @{Cell="I22"; Value=1.037057860165058},
    @{Cell="J22"; Value=1.040999866788773},
    @{Cell="K22"; Value=1.040262734866942},
    @{Cell="L22"; Value=1.046249224886867},
    @{Cell="M22"; Value=1.05511392349133},
    @{Cell="N22"; Value=1.017540546595247},
    @{Cell="B23"; Value=1.02},
    @{Cell="C23"; Value=1.034525287689766},
    @{Cell="D23"; Value=1.036824941978105},
    @{Cell="E23"; Value=1.042919542637768},
    @{Cell="F23"; Value=1.05191061125015},
    @{Cell="I23"; Value=1.037148903654406},
    @{Cell="J23"; Value=1.041236609624494},
    @{Cell="K23"; Value=1.040462649985046},
    @{Cell="L23"; Value=1.046534348650875},
    @{Cell="M23"; Value=1.055492166888793},
    @{Cell="N23"; Value=1.017620672666158},
    @{Cell="B24"; Value=1.02},
    @{Cell="C24"; Value=1.035978101481692},
    @{Cell="D24"; Value=1.037885043719293},
    @{Cell="E24"; Value=1.04431554188864},
    @{Cell="F24"; Value=1.053672946558295},
    @{Cell="I24"; Value=1.037504918078311},
    @{Cell="J24"; Value=1.04216780065451},
    @{Cell="K24"; Value=1.041248540136579},
    @{Cell="L24"; Value=1.047656839100386},
    @{Cell="M24"; Value=1.056982472233595},
    @{Cell="N24"; Value=1.017935656889111},
    @{Cell="B25"; Value=1.02},
    @{Cell="C25"; Value=1.037665032851268},
    @{Cell="D25"; Value=1.039115961574265},
    @{Cell="E25"; Value=1.045938376863591},
    @{Cell="F25"; Value=1.055723154746014},
    @{Cell="I25"; Value=1.037913039685695},
    @{Cell="J25"; Value=1.043246632796047},
    @{Cell="K25"; Value=1.042158106069105},
    @{Cell="L25"; Value=1.048959381315135},
    @{Cell="M25"; Value=1.058714360004478},
    @{Cell="N25"; Value=1.018300207693377}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
